# Commit message: "updates to all docs"
# This script:
#  1. Clears the stray "MISSING" legend label and the leftover -1/1 grid
#     underneath it (B19:E21), including their red-highlight formatting.
#  2. Resets row 21's row-level formatting (height/thick-bottom-border flag)
#     now that nothing in the row needs it.
#  3. Restores rows 12 and 16 (columns A:F) to the normal (non-highlighted)
#     look shared by every other data row, by pulling the formatting from a
#     neighboring, correctly-formatted row (row 13) without touching the
#     values/formulas already in row 12/16.
#  4. Updates the saved cursor/selection position in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the old "MISSING" legend entry and its stray data grid ---
$ws.Range("B19:E21").ClearContents()
$ws.Range("B19:E21").ClearFormats()
$ws.Range("B19:E21").Style = "Normal"

# --- 2. Row 21 no longer needs special height / thick-bottom-border ---
$ws.Rows.Item(21).AutoFit()

# --- 3. Re-normalize rows 12 and 16 (A:F) to match the rest of the table ---
$ws.Range("A13:F13").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$ws.Range("A16:F16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. Update the sheet view's saved scroll position / selection ---
$ws.Range("F19").Select()
$excel.ActiveWindow.ScrollRow = 4
